$wb = $excel.ActiveWorkbook

# --- classNumberOfLines sheet -------------------------------------------
# "Number of Lines" for pl.piomin.stock.domain.Product: 0 -> 39
# "Number of Lines" for pl.piomin.stock.repository.ProductRepository: 0 -> 1
$wsClass = $wb.Worksheets.Item("classNumberOfLines")
$wsClass.Range("B2").NumberFormat = "@"
$wsClass.Range("B4").NumberFormat = "@"
$wsClass.Range("B2").Value = "39"
$wsClass.Range("B4").Value = "1"

# --- methodNumberOfLines sheet ------------------------------------------
# Fix the "Number of Lines" (column C) for every method row so it reflects
# the real size of each method, instead of the placeholder 0.
$wsMethod = $wb.Worksheets.Item("methodNumberOfLines")
$wsMethod.Range("C2:C12").NumberFormat = "@"
$wsMethod.Range("C14").NumberFormat = "@"
$wsMethod.Range("C18").NumberFormat = "@"
$wsMethod.Range("C23").NumberFormat = "@"
$wsMethod.Range("C27").NumberFormat = "@"

$wsMethod.Range("C2").Value = "3"
$wsMethod.Range("C3").Value = "3"
$wsMethod.Range("C4").Value = "3"
$wsMethod.Range("C5").Value = "3"
$wsMethod.Range("C6").Value = "3"
$wsMethod.Range("C7").Value = "3"
$wsMethod.Range("C8").Value = "3"
$wsMethod.Range("C9").Value = "3"
$wsMethod.Range("C10").Value = "3"
$wsMethod.Range("C11").Value = "2"
$wsMethod.Range("C12").Value = "6"
$wsMethod.Range("C14").Value = "1"
$wsMethod.Range("C18").Value = "1"
$wsMethod.Range("C23").Value = "1"
$wsMethod.Range("C27").Value = "1"
